# Adding XML Dom File - checked in by Rudy Xiong
#
# Updates the PSO Public Machines list:
#  - fixes XP_15's PC Name (was mistakenly the same IP as XP_18)
#  - removes the old mailto: hyperlinks (and their blue/underlined
#    "Hyperlink" styling) that used to sit on the Password column
#  - appends two newly available machines: Win7_56 and XP_18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Strip the mailto: hyperlinks off the Password column (E2:E3) ---
$ws.Hyperlinks.Delete()

# ...and put those two cells back to plain, unformatted text instead of
# the underlined/blue "Hyperlink" look.
$ws.Range("E2").ClearFormats()
$ws.Range("E3").ClearFormats()

# That named cell style is now unused anywhere in the workbook - drop it.
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "超链接") {
        $s.Delete()
    }
}

# --- XP_15's PC Name was a copy/paste duplicate of XP_18's IP; fix it ---
$ws.Range("B2").Value = "10.224.168.15"

# --- Two more machines are now available; append them as rows 4 and 5 ---
$ws.Range("A4").Value = "Win7_56"
$ws.Range("B4").Value = "10.224.168.56"
$ws.Range("C4").Value = "Cisco"
$ws.Range("D4").Value = "cisco"
$ws.Range("E4").Value = "Pass"
$ws.Range("F4").Value = "Y"

$ws.Range("A5").Value = "XP_18"
$ws.Range("B5").Value = "10.224.168.18"
$ws.Range("C5").Value = "WebexHZ"
$ws.Range("D5").Value = "cisco"
$ws.Range("E5").Value = "pass"
$ws.Range("F5").Value = "N"
$ws.Range("G5").Value = "Evelyn Yao"

# --- Leave the selection where the new data ends, same as the author did ---
$ws.Range("G5").Select()
